$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-23 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2026-02-24 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("469÷2=234, 1", $false, $false, $false, $false, $false, $true, 1, $false, "492÷9=54, 6", 2) | Out-Null
$d.Content.Find.Execute("865÷9=96, 1", $false, $false, $false, $false, $false, $true, 1, $false, "493÷7=70, 3", 2) | Out-Null
$d.Content.Find.Execute("813÷6=135, 3", $false, $false, $false, $false, $false, $true, 1, $false, "791÷9=87, 8", 2) | Out-Null
$d.Content.Find.Execute("282÷3=94, 0", $false, $false, $false, $false, $false, $true, 1, $false, "815÷8=101, 7", 2) | Out-Null
$d.Content.Find.Execute("980÷3=326, 2", $false, $false, $false, $false, $false, $true, 1, $false, "230÷6=38, 2", 2) | Out-Null
$d.Content.Find.Execute("887÷3=295, 2", $false, $false, $false, $false, $false, $true, 1, $false, "463÷9=51, 4", 2) | Out-Null
$d.Content.Find.Execute("804÷3=268, 0", $false, $false, $false, $false, $false, $true, 1, $false, "729÷8=91, 1", 2) | Out-Null
$d.Content.Find.Execute("160÷9=17, 7", $false, $false, $false, $false, $false, $true, 1, $false, "966÷4=241, 2", 2) | Out-Null
$d.Content.Find.Execute("669÷6=111, 3", $false, $false, $false, $false, $false, $true, 1, $false, "120÷4=30, 0", 2) | Out-Null
$d.Content.Find.Execute("290÷6=48, 2", $false, $false, $false, $false, $false, $true, 1, $false, "383÷9=42, 5", 2) | Out-Null
$d.Content.Find.Execute("716÷9=79, 5", $false, $false, $false, $false, $false, $true, 1, $false, "582÷9=64, 6", 2) | Out-Null
$d.Content.Find.Execute("694÷2=347, 0", $false, $false, $false, $false, $false, $true, 1, $false, "494÷9=54, 8", 2) | Out-Null
$d.Content.Find.Execute("958÷7=136, 6", $false, $false, $false, $false, $false, $true, 1, $false, "266÷6=44, 2", 2) | Out-Null
$d.Content.Find.Execute("587÷3=195, 2", $false, $false, $false, $false, $false, $true, 1, $false, "718÷9=79, 7", 2) | Out-Null
$d.Content.Find.Execute("535÷4=133, 3", $false, $false, $false, $false, $false, $true, 1, $false, "759÷3=253, 0", 2) | Out-Null
$d.Content.Find.Execute("298÷4=74, 2", $false, $false, $false, $false, $false, $true, 1, $false, "562÷3=187, 1", 2) | Out-Null
$d.Content.Find.Execute("126÷8=15, 6", $false, $false, $false, $false, $false, $true, 1, $false, "977÷7=139, 4", 2) | Out-Null
$d.Content.Find.Execute("519÷6=86, 3", $false, $false, $false, $false, $false, $true, 1, $false, "238÷2=119, 0", 2) | Out-Null
$d.Content.Find.Execute("262÷9=29, 1", $false, $false, $false, $false, $false, $true, 1, $false, "769÷4=192, 1", 2) | Out-Null
$d.Content.Find.Execute("209÷3=69, 2", $false, $false, $false, $false, $false, $true, 1, $false, "865÷3=288, 1", 2) | Out-Null
$d.Content.Find.Execute("606÷8=75, 6", $false, $false, $false, $false, $false, $true, 1, $false, "275÷4=68, 3", 2) | Out-Null
$d.Content.Find.Execute("155÷7=22, 1", $false, $false, $false, $false, $false, $true, 1, $false, "650÷2=325, 0", 2) | Out-Null
$d.Content.Find.Execute("298÷9=33, 1", $false, $false, $false, $false, $false, $true, 1, $false, "572÷6=95, 2", 2) | Out-Null
$d.Content.Find.Execute("440÷7=62, 6", $false, $false, $false, $false, $false, $true, 1, $false, "362÷8=45, 2", 2) | Out-Null
$d.Content.Find.Execute("931÷6=155, 1", $false, $false, $false, $false, $false, $true, 1, $false, "173÷7=24, 5", 2) | Out-Null
